# Generate Report for Handback
#
# For each localized-language source file that has come back from
# localization (rows 2 & 3 of the "zh-cn" / "de-de" sheets), the handback
# report now records:
#   - an updated Status ("Ready for handoff" -> "Handed back: in sync with en-US")
#   - the "Latest Target File" (column E) - same file as the source (column A),
#     hyperlinked the same way
#   - the "Latest Handback File" (column F) - same file as the latest handoff
#     (column C), hyperlinked the same way
#   - the "Latest Handback DateTime" (column G) - the timestamp the handback
#     was produced
#
# The Overview sheet shows the same Status text for each file/language pair,
# so it is refreshed too.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet - refresh the Status column for both languages
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

$zhcn.Hyperlinks.Add(
    $zhcn.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0c015fa9d6d55f3a084a3af2c4b4e9db8a6e5cb7/e2e/430d3873-a87e-4a4a-9e4b-b5133148f1c2.md",
    "",
    "",
    "430d3873-a87e-4a4a-9e4b-b5133148f1c2.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60c4785db8864f4f2f6445ab5caf9c58ba16b0ff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/430d3873-a87e-4a4a-9e4b-b5133148f1c2.e96a723b2438091fe7701b4da476f71bd7e1b460.zh-cn.xlf",
    "",
    "",
    "430d3873-a87e-4a4a-9e4b-b5133148f1c2.e96a723b2438091fe7701b4da476f71bd7e1b460.zh-cn.xlf"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0c015fa9d6d55f3a084a3af2c4b4e9db8a6e5cb7/e2e/f87603a8-2281-4998-835f-d1e898eb12cd.md",
    "",
    "",
    "f87603a8-2281-4998-835f-d1e898eb12cd.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60c4785db8864f4f2f6445ab5caf9c58ba16b0ff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f87603a8-2281-4998-835f-d1e898eb12cd.fda8150841d30312ef1c6b7f5cd4e1bf2ead05ac.zh-cn.xlf",
    "",
    "",
    "f87603a8-2281-4998-835f-d1e898eb12cd.fda8150841d30312ef1c6b7f5cd4e1bf2ead05ac.zh-cn.xlf"
)

$zhcn.Range("G2").Value = "2016-03-08 21:18:01"
$zhcn.Range("G3").Value = "2016-03-08 21:18:01"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

$dede.Hyperlinks.Add(
    $dede.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0c015fa9d6d55f3a084a3af2c4b4e9db8a6e5cb7/e2e/430d3873-a87e-4a4a-9e4b-b5133148f1c2.md",
    "",
    "",
    "430d3873-a87e-4a4a-9e4b-b5133148f1c2.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c50d5b72448503badc2242994b59733bb94faef7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/430d3873-a87e-4a4a-9e4b-b5133148f1c2.e96a723b2438091fe7701b4da476f71bd7e1b460.de-de.xlf",
    "",
    "",
    "430d3873-a87e-4a4a-9e4b-b5133148f1c2.e96a723b2438091fe7701b4da476f71bd7e1b460.de-de.xlf"
)
$dede.Hyperlinks.Add(
    $dede.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0c015fa9d6d55f3a084a3af2c4b4e9db8a6e5cb7/e2e/f87603a8-2281-4998-835f-d1e898eb12cd.md",
    "",
    "",
    "f87603a8-2281-4998-835f-d1e898eb12cd.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c50d5b72448503badc2242994b59733bb94faef7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f87603a8-2281-4998-835f-d1e898eb12cd.fda8150841d30312ef1c6b7f5cd4e1bf2ead05ac.de-de.xlf",
    "",
    "",
    "f87603a8-2281-4998-835f-d1e898eb12cd.fda8150841d30312ef1c6b7f5cd4e1bf2ead05ac.de-de.xlf"
)

$dede.Range("G2").Value = "2016-03-08 21:18:25"
$dede.Range("G3").Value = "2016-03-08 21:18:25"

Write-Output "Handback report generated"
